$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "E19", "E20", "E21", "D22", "E22", "D23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D40", "E40", "D41", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50")
foreach ($c in $cells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "245.92"
$ws.Range("E2").Value = "-0.23%"
$ws.Range("D3").Value = "29.75"
$ws.Range("E3").Value = "-0.06%"
$ws.Range("D4").Value = "5.336"
$ws.Range("E4").Value = "2.84%"
$ws.Range("D5").Value = "0.05753"
$ws.Range("E5").Value = "0.88%"
$ws.Range("D6").Value = "6.651"
$ws.Range("E6").Value = "0.90%"
$ws.Range("D7").Value = "3.221"
$ws.Range("E7").Value = "5.64%"
$ws.Range("D8").Value = "0.8580"
$ws.Range("E8").Value = "0.02%"
$ws.Range("D9").Value = "0.8559"
$ws.Range("E9").Value = "-2.25%"
$ws.Range("D10").Value = "0.1382"
$ws.Range("E10").Value = "1.05%"
$ws.Range("D11").Value = "0.07099"
$ws.Range("E11").Value = "0.12%"
$ws.Range("D12").Value = "0.03193"
$ws.Range("E12").Value = "11.57%"
$ws.Range("D13").Value = "0.09338"
$ws.Range("D14").Value = "0.001523"
$ws.Range("E14").Value = "-1.55%"
$ws.Range("D15").Value = "0.0005934"
$ws.Range("E15").Value = "-1.62%"
$ws.Range("D16").Value = "0.005953"
$ws.Range("E16").Value = "-1.81%"
$ws.Range("E17").Value = "0.97%"
$ws.Range("D18").Value = "2.196"
$ws.Range("E18").Value = "-2.63%"
$ws.Range("E19").Value = "0.45%"
$ws.Range("E20").Value = "2.47%"
$ws.Range("E21").Value = "0.25%"
$ws.Range("D22").Value = "3.479"
$ws.Range("E22").Value = "19.71%"
$ws.Range("D23").Value = "0.04127"
$ws.Range("D24").Value = "0.1409"
$ws.Range("E24").Value = "2.06%"
$ws.Range("D25").Value = "0.001222"
$ws.Range("E25").Value = "0.27%"
$ws.Range("D26").Value = "0.004176"
$ws.Range("E26").Value = "-17.89%"
$ws.Range("D27").Value = "0.0001199"
$ws.Range("E27").Value = "-0.90%"
$ws.Range("D28").Value = "0.0001448"
$ws.Range("E28").Value = "-25.32%"
$ws.Range("D40").Value = "0.03773"
$ws.Range("E40").Value = "0.67%"
$ws.Range("D41").Value = "0.1073"
$ws.Range("E41").Value = "0.17%"
$ws.Range("D42").Value = "0.002414"
$ws.Range("E42").Value = "14.99%"
$ws.Range("E43").Value = "-48.10%"
$ws.Range("D44").Value = "0.009200"
$ws.Range("E44").Value = "-2.21%"
$ws.Range("D45").Value = "0.00005264"
$ws.Range("E45").Value = "2.98%"
$ws.Range("E46").Value = "-0.06%"
$ws.Range("D47").Value = "0.08082"
$ws.Range("E47").Value = "13.78%"
$ws.Range("D48").Value = "0.002203"
$ws.Range("E48").Value = "-17.52%"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").Value = "-0.06%"
